# 2019 12 06 21:25 inclusao das paginas de consulta
# Updates the "Periodo" duty-roster assignments in the December schedule
# (a rotation/shuffle of the existing names across several rows) and
# refreshes the "Data de geracao" generation timestamp in C35.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Jane"
$ws.Range("D4").Value = "Lucia"
$ws.Range("F4").Value = "Alana"
$ws.Range("G4").Value = "Graca"
$ws.Range("J4").Value = "Antonio"
$ws.Range("K4").Value = "Vinicius"
$ws.Range("L4").Value = "Vanda"
$ws.Range("M4").Value = "EMPTY"
$ws.Range("D5").Value = "Vanda"
$ws.Range("J5").Value = "Jessica Silva"
$ws.Range("K5").Value = "Daiana"
$ws.Range("C6").Value = "Rodolfo Dias"
$ws.Range("D6").Value = "Antonio"
$ws.Range("F6").Value = "Lindoia"
$ws.Range("G6").Value = "Edith"
$ws.Range("K6").Value = "Patricia Dias"
$ws.Range("D7").Value = "Daniel"
$ws.Range("C8").Value = "Robson"
$ws.Range("D8").Value = "Antonio"
$ws.Range("G8").Value = "Lucia"
$ws.Range("J8").Value = "Vinicius"
$ws.Range("K8").Value = "EMPTY"
$ws.Range("L8").Value = "Icaro"
$ws.Range("M8").Value = "Rodolfo Dias"
$ws.Range("K9").Value = "Beth"
$ws.Range("M9").Value = "Alex"
$ws.Range("C11").Value = "Helaine Camilo"
$ws.Range("D11").Value = "Aline Lima"
$ws.Range("F11").Value = "Lurdes"
$ws.Range("J11").Value = "Valquiria"
$ws.Range("K11").Value = "Patricia Dias"
$ws.Range("L11").Value = "Dario"
$ws.Range("M11").Value = "Amintas"
$ws.Range("D12").Value = "Cida"
$ws.Range("G12").Value = "Patricia Rodrigues"
$ws.Range("J12").Value = "Icaro"
$ws.Range("K12").Value = "Rodolfo Dias"
$ws.Range("D13").Value = "Valquiria"
$ws.Range("F13").Value = "Lurdes"
$ws.Range("J13").Value = "Vanda"
$ws.Range("K13").Value = "EMPTY"
$ws.Range("D14").Value = "Lindoia"
$ws.Range("C15").Value = "Isabele"
$ws.Range("D15").Value = "Daiana"
$ws.Range("G15").Value = "Patricia Rodrigues"
$ws.Range("J15").Value = "Beth"
$ws.Range("K15").Value = "Aline Lima"
$ws.Range("M15").Value = "EMPTY"
$ws.Range("K16").Value = "Carlos Eduardo"
$ws.Range("M16").Value = "Douglas Oliveira"
$ws.Range("C18").Value = "Lucia"
$ws.Range("D18").Value = "Jane"
$ws.Range("F18").Value = "Graca"
$ws.Range("G18").Value = "Alana"
$ws.Range("J18").Value = "Keila"
$ws.Range("K18").Value = "Eliane"
$ws.Range("L18").Value = "Marcio"
$ws.Range("D19").Value = "Keila"
$ws.Range("G19").Value = "Alana"
$ws.Range("J19").Value = "Jessica Silva"
$ws.Range("K19").Value = "Daiana"
$ws.Range("D20").Value = "Lindoia"
$ws.Range("F20").Value = "Rodolfo Dias"
$ws.Range("G20").Value = "Antonio"
$ws.Range("J20").Value = "Lucia"
$ws.Range("K20").Value = "Patricia Dias"
$ws.Range("D21").Value = "Keila"
$ws.Range("C22").Value = "Antonio"
$ws.Range("D22").Value = "Daniel"
$ws.Range("G22").Value = "Valquiria"
$ws.Range("J22").Value = "Keila"
$ws.Range("K22").Value = "Lurdes"
$ws.Range("L22").Value = "Rodolfo Dias"
$ws.Range("M22").Value = "Icaro"
$ws.Range("K23").Value = "Marcio"
$ws.Range("M23").Value = "Icaro"
$ws.Range("C25").Value = "Helaine Camilo"
$ws.Range("D25").Value = "Aline Lima"
$ws.Range("F25").Value = "Lurdes"
$ws.Range("J25").Value = "Patricia Dias"
$ws.Range("K25").Value = "Valquiria"
$ws.Range("L25").Value = "Vanda"
$ws.Range("M25").Value = "EMPTY"
$ws.Range("D26").Value = "Vanda"
$ws.Range("J26").Value = "Icaro"
$ws.Range("K26").Value = "Rodolfo Dias"
$ws.Range("C27").Value = "Rodolfo Dias"
$ws.Range("D27").Value = "Aline Lima"
$ws.Range("F27").Value = "Valquiria"
$ws.Range("G27").Value = "Lurdes"
$ws.Range("J27").Value = "Vanda"
$ws.Range("C29").Value = "Isabele"
$ws.Range("D29").Value = "Eliane"
$ws.Range("G29").Value = "Patricia Rodrigues"
$ws.Range("J29").Value = "Aline Lima"
$ws.Range("K29").Value = "Lucia"
$ws.Range("M29").Value = "Vinicius"
$ws.Range("K30").Value = "Alex"
$ws.Range("M30").Value = "Geronimo"
$ws.Range("C32").Value = "Alana"
$ws.Range("D32").Value = "Lucia"
$ws.Range("F32").Value = "Graca"
$ws.Range("J32").Value = "Keila"
$ws.Range("K32").Value = "Eliane"
$ws.Range("L32").Value = "Antonio"
$ws.Range("M32").Value = "Amintas"
$ws.Range("D33").Value = "Eliane"
$ws.Range("G33").Value = "Daiana"
$ws.Range("J33").Value = "Jessica Silva"
$ws.Range("K33").Value = "Alana"

# Update generation date/time (stored as Excel serial date-time number)
$ws.Range("C35").Value = 43807.75785267639
